$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the workpackage name and log file path strings to reflect the
# renamed process ("RPA Bank Change Letter" -> "RPA Bank Account Change Letter").
$ws.Range("B19").Value = "RPA Bank Account Change Letter"
$ws.Range("B18").Value = "C:\Users\{0}\Desktop\Bank Account Change Letter Logs_{1}.xlsx"

# Update the view so the active cell is B19 and the sheet is scrolled back
# to the top (no frozen/topLeft offset).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B19").Select()
